# Populate buyer order rows + style the header and data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data --------------------------------------------------------------
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Apple Juice: 5pcs., Manggo Juice: 1pcs., Guyabano Juice: 1pcs."
$ws.Range("C2").Value = "PHP 140.00"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Guyabano Juice: 2pcs., Ice Mixed Guyabano with Teal: 4pcs., Manggo Juice: 1pcs."
$ws.Range("C3").Value = "PHP 180.00"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Apple Juice: 3pcs."
$ws.Range("C4").Value = "PHP 60.00"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Ice Mixed Guyabano with Teal: 3pcs., Manggo Juice: 3pcs., Guyabano Juice: 1pcs."
$ws.Range("C5").Value = "PHP 170.00"

# --- Header formatting (row 1): bold white text on dark slate gray -----
$header = $ws.Range("A1:C1")
$header.Font.Bold = $true
$header.Font.Color = 16777215        # white
$header.Interior.Color = 5197615     # dark slate gray (2F4F4F)

# --- Column A (Buyer ID) formatting: bold, light blue, left/center -----
$colA = $ws.Range("A2:A5")
$colA.Font.Bold = $true
$colA.Interior.Color = 15128749      # light blue (ADD8E6)
$colA.HorizontalAlignment = -4131    # xlHAlignLeft
$colA.VerticalAlignment = -4108      # xlVAlignCenter

# --- Column B (Items Ordered) formatting: bold, light gray, left/center
$colB = $ws.Range("B2:B5")
$colB.Font.Bold = $true
$colB.Interior.Color = 13882323      # light gray (D3D3D3)
$colB.HorizontalAlignment = -4131    # xlHAlignLeft
$colB.VerticalAlignment = -4108      # xlVAlignCenter

# --- Column C (Total) formatting: bold, light pink/red, center/center --
$colC = $ws.Range("C2:C5")
$colC.Font.Bold = $true
$colC.Interior.Color = 13421823      # light pink (FFCCCC)
$colC.HorizontalAlignment = -4108    # xlHAlignCenter
$colC.VerticalAlignment = -4108      # xlVAlignCenter
